$wb = $excel.ActiveWorkbook

# --- Production_Rate sheet ---
$wsProd = $wb.Worksheets.Item("Production_Rate")
$wsProd.Range("B2").Value = 20.25889121338912
$wsProd.Range("C2").Value = 486.2133891213389
$wsProd.Range("B3").Value = 42.29901960784314
$wsProd.Range("C3").Value = 1015.176470588235

# --- Hourly_Arrivals sheet ---
$wsArr = $wb.Worksheets.Item("Hourly_Arrivals")
$wsArr.Range("B2").Value = 0.004184100418410041
$wsArr.Range("B3").Value = 0.004184100418410041
$wsArr.Range("B4").Value = 2.456066945606695
$wsArr.Range("B5").Value = 4.138075313807532
$wsArr.Range("B6").Value = 3.740585774058578
$wsArr.Range("B7").Value = 3.912133891213389
$wsArr.Range("B8").Value = 4.133891213389122
$wsArr.Range("B9").Value = 4.221757322175733
$wsArr.Range("B10").Value = 3.94142259414226
$wsArr.Range("B11").Value = 3.841004184100418
$wsArr.Range("B12").Value = 3.907949790794979
$wsArr.Range("B13").Value = 3.644351464435146
$wsArr.Range("B14").Value = 3.510460251046025
$wsArr.Range("B15").Value = 3.01255230125523
$wsArr.Range("B16").Value = 2.778242677824268
$wsArr.Range("B17").Value = 3.125523012552301
$wsArr.Range("B18").Value = 1.640167364016736
$wsArr.Range("B19").Value = 1.225941422594142
$wsArr.Range("B20").Value = 0.7949790794979079
$wsArr.Range("B21").Value = 0.4476987447698745
